$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.670.10'
$ws.Range('E2').Value = '  +0.74%  '
$ws.Range('D3').Value = '2.509.37'
$ws.Range('E3').Value = '  -1.41%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '590.66'
$ws.Range('E5').Value = '  -0.20%  '
$ws.Range('D6').Value = '172.45'
$ws.Range('E6').Value = '  -0.79%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').Value = '0.524'
$ws.Range('E8').Value = '  -0.94%  '
$ws.Range('D9').Value = '2.504.00'
$ws.Range('E9').Value = '  -1.63%  '
$ws.Range('E10').Value = '  -0.46%  '
$ws.Range('E11').Value = '  +1.20%  '
$ws.Range('D12').Value = '5.11'
$ws.Range('E12').Value = '  -0.79%  '
$ws.Range('D13').Value = '0.340'
$ws.Range('E13').Value = '  -1.75%  '
$ws.Range('D14').Value = '26.35'
$ws.Range('E14').Value = '  -1.76%  '
$ws.Range('D15').Value = '2.944.67'
$ws.Range('E15').Value = '  -2.16%  '
$ws.Range('D16').Value = '0.0000177'
$ws.Range('E16').Value = '  -0.10%  '
$ws.Range('D17').Value = '67.616.64'
$ws.Range('E17').Value = '  +1.12%  '
$ws.Range('D18').Value = '2.470.78'
$ws.Range('E18').Value = '  -2.44%  '
$ws.Range('D19').Value = '11.75'
$ws.Range('E19').Value = '  +3.27%  '
$ws.Range('D20').Value = '7.85'
$ws.Range('E20').Value = '  -2.94%  '
$ws.Range('D21').Value = '367.13'
$ws.Range('E21').Value = '  +3.35%  '
$ws.Range('D22').Value = '4.16'
$ws.Range('E22').Value = '  -0.91%  '
$ws.Range('D23').Value = '4.57'
$ws.Range('E23').Value = '  -1.08%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').Value = '71.40'
$ws.Range('E24').Value = '  +2.12%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('D26').Value = '1.93'
$ws.Range('E26').Value = '  -3.44%  '
$ws.Range('E27').Value = '  -1.72%  '
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.13%  '
$ws.Range('D29').Value = '2.601.94'
$ws.Range('E29').Value = '  -2.88%  '
$ws.Range('D30').Value = '0.0₃0960'
$ws.Range('E30').Value = '  -2.91%  '
$ws.Range('D31').Value = '8.34'
$ws.Range('E31').Value = '  +1.83%  '
$ws.Range('D32').Value = '533.58'
$ws.Range('E32').Value = '  -0.39%  '
$ws.Range('D33').Value = '1.32'
$ws.Range('E33').Value = '  -2.43%  '
$ws.Range('D34').Value = '1.87'
$ws.Range('E34').Value = '  +0.79%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').Value = '0.128'
$ws.Range('E35').Value = '  -3.26%  '
$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D36').Value = '0.998'
$ws.Range('E36').Value = '  -0.20%  '
$ws.Range('D37').Value = '158.48'
$ws.Range('E37').Value = '  +0.92%  '
$ws.Range('D38').Value = '1.42'
$ws.Range('E38').Value = '  -2.93%  '
$ws.Range('D39').Value = '18.84'
$ws.Range('E39').Value = '  +0.90%  '
$ws.Range('D40').Value = '18.63'
$ws.Range('E40').Value = '  +0.93%  '
$ws.Range('D41').Value = '0.349'
$ws.Range('E41').Value = '  -2.01%  '
$ws.Range('D42').Value = '1.78'
$ws.Range('E42').Value = '  -1.02%  '
$ws.Range('D43').Value = '5.12'
$ws.Range('E43').Value = '  -0.48%  '
$ws.Range('E44').Value = '  +0.10%  '
$ws.Range('D45').Value = '2.46'
$ws.Range('E45').Value = '  -1.94%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '145.86'
$ws.Range('E46').Value = '  -2.51%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.0₆0279'
$ws.Range('E47').Value = '  -0.41%  '
$ws.Range('D48').Value = '3.69'
$ws.Range('E48').Value = '  -0.31%  '
$ws.Range('D49').Value = '0.549'
$ws.Range('E49').Value = '  -2.11%  '
$ws.Range('D50').Value = '1.71'
$ws.Range('E50').Value = '  +0.47%  '
$ws.Range('D51').Value = '0.0749'
$ws.Range('E51').Value = '  -1.50%  '
